# Update attendee/sales figures (column F) on the "展览" (sheet1) and
# "全部类型" (sheet4, the combined view) sheets, plus the single figure
# on "演出" (sheet2), to match the refreshed data pulled at commit 456a3b4.
# "本地生活" (sheet3) has no data rows, so it needs no edits.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")      # 展览
$ws2 = $wb.Worksheets.Item("演出")      # 演出
$ws4 = $wb.Worksheets.Item("全部类型")  # 全部类型

# --- 展览 (exhibitions) sheet ---
$ws1.Range("F4").Value = 300
$ws1.Range("F7").Value = 775
$ws1.Range("F8").Value = 75
$ws1.Range("F9").Value = 193
$ws1.Range("F11").Value = 473
$ws1.Range("F12").Value = 1423
$ws1.Range("F14").Value = 128
$ws1.Range("F15").Value = 289
$ws1.Range("F17").Value = 101
$ws1.Range("F18").Value = 678
$ws1.Range("F20").Value = 43
$ws1.Range("F21").Value = 268
$ws1.Range("F23").Value = 6125
$ws1.Range("F25").Value = 129
$ws1.Range("F26").Value = 124
$ws1.Range("F28").Value = 14867
$ws1.Range("F29").Value = 1476
$ws1.Range("F30").Value = 255
$ws1.Range("F31").Value = 116
$ws1.Range("F32").Value = 93
$ws1.Range("F33").Value = 10852
$ws1.Range("F34").Value = 682
$ws1.Range("F35").Value = 4247
$ws1.Range("F36").Value = 183
$ws1.Range("F37").Value = 367

# --- 演出 (performances) sheet ---
$ws2.Range("F2").Value = 340

# --- 全部类型 (all types, combined) sheet ---
$ws4.Range("F4").Value = 300
$ws4.Range("F7").Value = 775
$ws4.Range("F8").Value = 75
$ws4.Range("F9").Value = 193
$ws4.Range("F11").Value = 473
$ws4.Range("F12").Value = 1423
$ws4.Range("F14").Value = 128
$ws4.Range("F15").Value = 289
$ws4.Range("F16").Value = 340
$ws4.Range("F18").Value = 101
$ws4.Range("F19").Value = 678
$ws4.Range("F22").Value = 43
$ws4.Range("F23").Value = 268
$ws4.Range("F26").Value = 6125
$ws4.Range("F28").Value = 129
$ws4.Range("F29").Value = 124
$ws4.Range("F31").Value = 14867
$ws4.Range("F32").Value = 1476
$ws4.Range("F33").Value = 255
$ws4.Range("F34").Value = 116
$ws4.Range("F35").Value = 93
$ws4.Range("F36").Value = 10852
$ws4.Range("F37").Value = 682
$ws4.Range("F38").Value = 4247
$ws4.Range("F39").Value = 183
$ws4.Range("F40").Value = 367

Write-Output "Done updating attendance figures."
